$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.593583941459656
$ws.Range("B1").Value = 2.314446449279785
$ws.Range("C1").Value = 2.662900686264038
$ws.Range("D1").Value = 3.484760999679565
$ws.Range("E1").Value = 0.5208554863929749
